# Fruta / hortaliza, semanal
# Insert a new weekly record at row 153 (pushing the existing rows 153-177
# down to 154-178) in the Zapallo / Feria Lagunitas de Puerto Montt sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 153..177 down one row, creating space for the new record.
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the latest week's data.
$ws.Range("A153").Value = 4
$ws.Range("B153").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C153").Value = "Los Lagos"
$ws.Range("D153").Value = 44474
$ws.Range("E153").Value = 10
$ws.Range("F153").Value = 100112045
$ws.Range("G153").Value = "Zapallo"
$ws.Range("H153").Value = "Paine"
$ws.Range("I153").Value = "1a (guarda)"
$ws.Range("J153").Value = 1200
$ws.Range("K153").Value = 580
$ws.Range("L153").Value = 600
$ws.Range("M153").Value = 590
$ws.Range("N153").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O153").Value = "Región Metropolitana"
$ws.Range("P153").Value = 590
$ws.Range("Q153").Value = 1
$ws.Range("R153").Value = "Hortaliza"
